# Add team-record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF -----------------------
# Copy the formatting of the last existing header cell (AC1, style "s=1":
# bold font + thin border + centered/top alignment) onto the three new
# header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-68): same team record for every player ---------
$wins = 62
$losses = 100
$ties = 0

for ($r = 2; $r -le 68; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-68"
